$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 39.288329
$ws.Range("H2").Value = 117.864987
$ws.Range("I2").Value = 0.632237668435316
$ws.Range("J2").Value = 0.632237668435316
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.970048
$ws.Range("N2").Value = 26.910144
$ws.Range("O2").Value = 0.487108783009476
$ws.Range("P2").Value = 0.4871087830094759
$ws.Range("Q2").Value = 352.418196969792
$ws.Range("R2").Value = 3171.763772728128
$ws.Range("S2").Value = 0.3079685212442754
$ws.Range("T2").Value = 0.3079685212442753
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 39.288329
$ws.Range("H3").Value = 117.864987
$ws.Range("I3").Value = 0.632237668435316
$ws.Range("J3").Value = 0.632237668435316
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.012070666666666
$ws.Range("N3").Value = 27.036212
$ws.Range("O3").Value = 0.489390778604016
$ws.Range("P3").Value = 0.489390778604016
$ws.Range("Q3").Value = 354.0691973232493
$ws.Range("R3").Value = 3186.622775909244
$ws.Range("S3").Value = 0.309411284818347
$ws.Range("T3").Value = 0.309411284818347
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 39.288329
$ws.Range("H4").Value = 117.864987
$ws.Range("I4").Value = 0.632237668435316
$ws.Range("J4").Value = 0.632237668435316
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.4327576666666667
$ws.Range("N4").Value = 1.298273
$ws.Range("O4").Value = 0.02350043838650813
$ws.Range("P4").Value = 0.02350043838650813
$ws.Range("Q4").Value = 17.00232558527233
$ws.Range("R4").Value = 153.020930267451
$ws.Range("S4").Value = 0.0148578623726937
$ws.Range("T4").Value = 0.0148578623726937
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.344283
$ws.Range("H5").Value = 58.032849
$ws.Range("I5").Value = 0.3112930657211948
$ws.Range("J5").Value = 0.3112930657211947
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.970048
$ws.Range("N5").Value = 26.910144
$ws.Range("O5").Value = 0.487108783009476
$ws.Range("P5").Value = 0.4871087830094759
$ws.Range("Q5").Value = 173.519147035584
$ws.Range("R5").Value = 1561.672323320256
$ws.Range("S5").Value = 0.15163358640274
$ws.Range("T5").Value = 0.1516335864027399
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.344283
$ws.Range("H6").Value = 58.032849
$ws.Range("I6").Value = 0.3112930657211948
$ws.Range("J6").Value = 0.3112930657211947
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.012070666666666
$ws.Range("N6").Value = 27.036212
$ws.Range("O6").Value = 0.489390778604016
$ws.Range("P6").Value = 0.489390778604016
$ws.Range("Q6").Value = 174.3320453919987
$ws.Range("R6").Value = 1568.988408527988
$ws.Range("S6").Value = 0.1523439558073266
$ws.Range("T6").Value = 0.1523439558073266
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.344283
$ws.Range("H7").Value = 58.032849
$ws.Range("I7").Value = 0.3112930657211948
$ws.Range("J7").Value = 0.3112930657211947
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.4327576666666667
$ws.Range("N7").Value = 1.298273
$ws.Range("O7").Value = 0.02350043838650813
$ws.Range("P7").Value = 0.02350043838650813
$ws.Range("Q7").Value = 8.371386774419667
$ws.Range("R7").Value = 75.342480969777
$ws.Range("S7").Value = 0.007315523511128162
$ws.Range("T7").Value = 0.007315523511128161
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.509096666666667
$ws.Range("H8").Value = 10.52729
$ws.Range("I8").Value = 0.05646926584348937
$ws.Range("J8").Value = 0.05646926584348937
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.970048
$ws.Range("N8").Value = 26.910144
$ws.Range("O8").Value = 0.487108783009476
$ws.Range("P8").Value = 0.4871087830094759
$ws.Range("Q8").Value = 31.47676553664
$ws.Range("R8").Value = 283.29088982976
$ws.Range("S8").Value = 0.02750667536246067
$ws.Range("T8").Value = 0.02750667536246067
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.509096666666667
$ws.Range("H9").Value = 10.52729
$ws.Range("I9").Value = 0.05646926584348937
$ws.Range("J9").Value = 0.05646926584348937
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 9.012070666666666
$ws.Range("N9").Value = 27.036212
$ws.Range("O9").Value = 0.489390778604016
$ws.Range("P9").Value = 0.489390778604016
$ws.Range("Q9").Value = 31.62422713616444
$ws.Range("R9").Value = 284.61804422548
$ws.Range("S9").Value = 0.02763553797834243
$ws.Range("T9").Value = 0.02763553797834243
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.509096666666667
$ws.Range("H10").Value = 10.52729
$ws.Range("I10").Value = 0.05646926584348937
$ws.Range("J10").Value = 0.05646926584348937
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.4327576666666667
$ws.Range("N10").Value = 1.298273
$ws.Range("O10").Value = 0.02350043838650813
$ws.Range("P10").Value = 0.02350043838650813
$ws.Range("Q10").Value = 1.518588485574444
$ws.Range("R10").Value = 13.66729637017
$ws.Range("S10").Value = 0.00132705250268627
$ws.Range("T10").Value = 0.00132705250268627
